# Add a new "2022-Q3" sheet (fund holdings detail) right after "总计",
# pushing "2022-Q2", "2022-Q1", "2021-Q4", "2021-Q1" one slot later, and
# update the "总计" summary sheet with the new quarter's totals.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Build the new "2022-Q3" worksheet.
#    "2022-Q1" already has the right shape (header + 7 data rows, 8 used
#    rows total) and all the right cell styling, so clone it (keeps the
#    header/border/bold formatting and the index-column style exactly),
#    move the clone to the 2nd tab position, rename it, then overwrite
#    every cell with the 2022-Q3 figures.
# ---------------------------------------------------------------------
$totalSheet = $wb.Worksheets.Item(1)
$templateSheet = $wb.Worksheets.Item("2022-Q1")
$templateSheet.Copy($null, $totalSheet)

$q3 = $wb.Worksheets.Item(2)
$q3.Name = "2022-Q3"

# Header row (unchanged wording, just re-asserted for clarity/safety).
$q3.Range("B1").Value = "基金代码"
$q3.Range("C1").Value = "基金名称"
$q3.Range("D1").Value = "基金规模"
$q3.Range("E1").Value = "股票总仓位"
$q3.Range("F1").Value = "仓位占比"
$q3.Range("G1").Value = "持有市值(亿元)"
$q3.Range("H1").Value = "仓位排名"

# Data rows. Columns D/E/F/G store numeric-looking figures as *text* in
# this workbook (matches every other quarter sheet), so a leading "'" is
# used to stop Excel from auto-coercing them (and fund codes in column B)
# into numbers / dropping leading zeros.
$q3.Range("A2").Value = 0
$q3.Range("B2").Value = "'011383"
$q3.Range("C2").Value = "富安达医药创新混合"
$q3.Range("D2").Value = "'1.43"
$q3.Range("E2").Value = "'83.09"
$q3.Range("F2").Value = "'4.64"
$q3.Range("G2").Value = "'0.0664"
$q3.Range("H2").Value = 7

$q3.Range("A3").Value = 1
$q3.Range("B3").Value = "'010746"
$q3.Range("C3").Value = "富安达长三角区域主题混合"
$q3.Range("D3").Value = "'0.98"
$q3.Range("E3").Value = "'88.86"
$q3.Range("F3").Value = "'4.51"
$q3.Range("G3").Value = "'0.0442"
$q3.Range("H3").Value = 7

$q3.Range("A4").Value = 2
$q3.Range("B4").Value = "'001861"
$q3.Range("C4").Value = "富安达健康人生灵活配置混合A"
$q3.Range("D4").Value = "'0.51"
$q3.Range("E4").Value = "'86.44"
$q3.Range("F4").Value = "'4.73"
$q3.Range("G4").Value = "'0.0241"
$q3.Range("H4").Value = 6

$q3.Range("A5").Value = 3
$q3.Range("B5").Value = "'014708"
$q3.Range("C5").Value = "天弘臻选健康混合A"
$q3.Range("D5").Value = "'0.58"
$q3.Range("E5").Value = "'85.49"
$q3.Range("F5").Value = "'3.51"
$q3.Range("G5").Value = "'0.0204"
$q3.Range("H5").Value = 8

$q3.Range("A6").Value = 4
$q3.Range("B6").Value = "'005293"
$q3.Range("C6").Value = "诺德新旺灵活配置混合"
$q3.Range("D6").Value = "'0.38"
$q3.Range("E6").Value = "'93.45"
$q3.Range("F6").Value = "'5.06"
$q3.Range("G6").Value = "'0.0192"
$q3.Range("H6").Value = 10

$q3.Range("A7").Value = 5
$q3.Range("B7").Value = "'014709"
$q3.Range("C7").Value = "天弘臻选健康混合C"
$q3.Range("D7").Value = "'0.08"
$q3.Range("E7").Value = "'85.49"
$q3.Range("F7").Value = "'3.51"
$q3.Range("G7").Value = "'0.0028"
$q3.Range("H7").Value = 8

$q3.Range("A8").Value = 6
$q3.Range("B8").Value = "'014470"
$q3.Range("C8").Value = "富安达健康人生灵活配置混合C"
$q3.Range("D8").Value = "'0.01"
$q3.Range("E8").Value = "'86.44"
$q3.Range("F8").Value = "'4.73"
$q3.Range("G8").Value = "'0.0005"
$q3.Range("H8").Value = 6

# Drop the "quote prefix" style that a leading "'" left behind on B/D/E/F/G
# so those cells come back to the same un-styled look as the rest of the
# (unmodified) quarter sheets.
$q3.Range("B2:B8").Style = "Normal"
$q3.Range("D2:G8").Style = "Normal"

# ---------------------------------------------------------------------
# 2) Update the "总计" (summary) sheet: insert the new 2022-Q3 row at the
#    top of the data, push the existing four quarters down one row, and
#    add a new row for 2021-Q1 at the bottom.
# ---------------------------------------------------------------------
$total = $wb.Worksheets.Item("总计")

# Give the new row 6 the same index-column styling (s="2") as the rows
# above it before writing its value.
$total.Range("A5").Copy()
$total.Range("A6").PasteSpecial(-4122)

$total.Range("A2").Value = 0
$total.Range("B2").Value = "2022-Q3"
$total.Range("C2").Value = 7
$total.Range("D2").Value = 0.18

$total.Range("A3").Value = 1
$total.Range("B3").Value = "2022-Q2"
$total.Range("C3").Value = 4
$total.Range("D3").Value = 0.65

$total.Range("A4").Value = 2
$total.Range("B4").Value = "2022-Q1"
$total.Range("C4").Value = 7
$total.Range("D4").Value = 1.22

$total.Range("A5").Value = 3
$total.Range("B5").Value = "2021-Q4"
$total.Range("C5").Value = 1
$total.Range("D5").Value = 0.02

$total.Range("A6").Value = 4
$total.Range("B6").Value = "2021-Q1"
$total.Range("C6").Value = 5
$total.Range("D6").Value = 0.33
